$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection back up to the top of the data; the frozen
# pane's top-left cell naturally follows the freeze boundary (row 3) once
# the sheet is saved.
$ws.Range("A2").Select()

# Data updates: column C holds the daily new-case counts that feed the
# cumulative running total in column B (shared formula). Updating C
# recalculates B (and the dependent J/K columns) automatically.
$ws.Range("C443").Value = 22
$ws.Range("C461").Value = 37
$ws.Range("C462").Value = 27
$ws.Range("C463").Value = 25
$ws.Range("C464").Value = 5

# Rows 464/465 gain explicit zero entries in L/M (previously blank). Those
# two columns are pre-formatted as Text ("@"), so writing a numeric 0
# straight into .Value would be stored as the text string "0" instead of
# a real number. Borrow a same-bordered General-format cell's formatting
# just long enough to write the number, then copy the original Text
# formatting back from an untouched sibling cell in the same column so
# the cell's look (and style index) ends up exactly as it started.
$ws.Range("D3").Copy()
$ws.Range("L464").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("L464").Value = 0
$ws.Range("L462").Copy()
$ws.Range("L464").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C462").Copy()
$ws.Range("M464").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M464").Value = 0
$ws.Range("M462").Copy()
$ws.Range("M464").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D3").Copy()
$ws.Range("L465").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("L465").Value = 0
$ws.Range("L462").Copy()
$ws.Range("L465").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C462").Copy()
$ws.Range("M465").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M465").Value = 0
$ws.Range("M462").Copy()
$ws.Range("M465").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
